$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999272596618
$ws.Range("A2").Value = 0.99564781667934743
$ws.Range("A3").Value = 0.97603340789236714
$ws.Range("A4").Value = 0.96742267913264024
$ws.Range("A5").Value = 0.95925901665258628
$ws.Range("A6").Value = 0.93940539984420446
$ws.Range("A7").Value = 0.93389163311969503
$ws.Range("A8").Value = 0.92637064078811338
$ws.Range("A9").Value = 0.93032888771307354
$ws.Range("A10").Value = 0.9352939468179573
$ws.Range("A11").Value = 0.93619362933502748
$ws.Range("A12").Value = 0.93341407638150709
$ws.Range("A13").Value = 0.92212793458537312
$ws.Range("A14").Value = 0.91796081451057288
$ws.Range("A15").Value = 0.91536945353654398
$ws.Range("A16").Value = 0.91286304919101358
$ws.Range("A17").Value = 0.90915523249956087
$ws.Range("A18").Value = 0.9080463348242892
$ws.Range("A19").Value = 0.99459463093536993
$ws.Range("A20").Value = 0.98747768899208932
$ws.Range("A21").Value = 0.98607922159196426
$ws.Range("A22").Value = 0.98481472670228576
$ws.Range("A23").Value = 0.96899803930389372
$ws.Range("A24").Value = 0.95597657763923083
$ws.Range("A25").Value = 0.94951948991357915
$ws.Range("A26").Value = 0.93574965967785928
$ws.Range("A27").Value = 0.93090276109271253
$ws.Range("A28").Value = 0.90942878835936347
$ws.Range("A29").Value = 0.89415706331173139
$ws.Range("A30").Value = 0.8875866137665509
$ws.Range("A31").Value = 0.87993284481851997
$ws.Range("A32").Value = 0.87825354326913296
$ws.Range("A33").Value = 0.8777335458339286
